$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held "8822123 - Roberta Veloso Garcia" (old row 13, with no
# label in column A) is removed entirely; everything below it shifts up by one row.
$ws.Rows("13").Delete()

# Row 10 (Objetivos:) now shows the "Docentes responsáveis" value instead of the
# old objectives paragraph.
$ws.Range("B10").Value = "8822123 - Roberta Veloso Garcia"
$ws.Range("C10").Value = "8822123 - Roberta Veloso Garcia"

# Row 13 (was old row 14, "Programa resumido:") now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (was old row 16, "Programa:") now shows the activation date value that
# already lives in B8/C8 ("Ativação:"). Copy it so the cell keeps its original
# text type/format instead of Excel auto-converting the "01/01/2018" string into
# a date serial number.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Row 18 (was old row 19, "Método:") now shows the docentes responsáveis value.
$ws.Range("B18").Value = "8822123 - Roberta Veloso Garcia"
$ws.Range("C18").Value = "8822123 - Roberta Veloso Garcia"

# Row 19 (was old row 20, "Critério:") now shows the old "Método" evaluation text.
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# Row 20 (was old row 21, "Norma de recuperação:") now shows the old "Critério" text.
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# Row 21 (was old row 22, "Bibliografia:") now shows the old "Norma de recuperação" text.
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
